# Applies the "Updated symbol list" edit (Fri Dec 23 19:37:46 UTC 2022) to Sheet1.
# - Column D ("Price") cells are updated with refreshed quotes. These values are
#   stored as text (e.g. "246.33", "0.0005890") so the NumberFormat is forced to
#   "Text" ("@") before assignment, and the style is reset back to "Normal"
#   afterwards so no stray numeric formatting is left behind.
# - Columns B ("Coin"), C ("Link") and E ("Volume(1h)") for rows 10-18 shift as the
#   ranking list churns (new coins enter / existing coins move), so those text
#   cells are simply overwritten with their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$Value
    )
    $range = $ws.Range($Address)
    $range.NumberFormat = "@"
    $range.Value = $Value
    $range.Style = "Normal"
}

# ----- Column D price updates (rows with no accompanying coin/link change) -----
Set-TextCell "D2"  "246.33"
Set-TextCell "D3"  "22.21"
Set-TextCell "D4"  "5.350"
Set-TextCell "D5"  "0.05862"
Set-TextCell "D7"  "6.379"
Set-TextCell "D8"  "0.8128"
Set-TextCell "D9"  "0.9977"

# ----- Rows 10-18: ranking list shuffled, coin/link/price/rank-label updated -----
Set-TextCell "B10" "WazirX"
Set-TextCell "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D10" "0.1424"
Set-TextCell "E10" "9WazirXWRX"

Set-TextCell "B11" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C11" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D11" "0.03637"
Set-TextCell "E11" "10LiechtensteinCryptoassetsExchangeLCX"

Set-TextCell "B12" "MandalaExchangeToken"
Set-TextCell "C12" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D12" "0.07359"
Set-TextCell "E12" "11MandalaExchangeTokenMDX"

Set-TextCell "B13" "BitrueCoin"
Set-TextCell "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D13" "0.03001"
Set-TextCell "E13" "12BitrueCoinBTR"

Set-TextCell "B14" "MCDex"
Set-TextCell "C14" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D14" "4.174"
Set-TextCell "E14" "13MCDexMCB"

Set-TextCell "B15" "BitMartToken"
Set-TextCell "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D15" "0.09395"
Set-TextCell "E15" "14BitMartTokenBMX"

Set-TextCell "B16" "BitForexToken"
Set-TextCell "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D16" "0.001596"
Set-TextCell "E16" "15BitForexTokenBF"

Set-TextCell "B17" "CoinExToken"
Set-TextCell "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell "D17" "0.04838"
Set-TextCell "E17" "16CoinExTokenCET"

Set-TextCell "B18" "One"
Set-TextCell "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D18" "0.0005889"
Set-TextCell "E18" "17OneONE"

# ----- Remaining column D price updates -----
Set-TextCell "D19" "0.005971"
Set-TextCell "D20" "0.004083"
Set-TextCell "D21" "0.0009895"
Set-TextCell "D22" "0.0001100"
Set-TextCell "D23" "3.689"
Set-TextCell "D24" "2.207"
Set-TextCell "D25" "0.3248"
Set-TextCell "D41" "0.006478"
Set-TextCell "D42" "0.1075"
Set-TextCell "D43" "0.002409"
Set-TextCell "D44" "0.005227"
Set-TextCell "D45" "0.00005657"
Set-TextCell "D48" "0.08273"
